$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (distance 1000)
$ws.Range("B2").Value = 101394.0
$ws.Range("C2").Value = 0.45
$ws.Range("E2").Value = 57.214067598674475
$ws.Range("H2").Value = 0.8621410668963386
$ws.Range("I2").Value = 56.6996
$ws.Range("K2").Value = 57.7947

# Row 3 (distance 1200)
$ws.Range("B3").Value = 118818.0
$ws.Range("C3").Value = 0.422
$ws.Range("E3").Value = 70.82140614048376
$ws.Range("H3").Value = 1.1549534660891005
$ws.Range("I3").Value = 70.13782499999999
$ws.Range("J3").Value = 70.91
$ws.Range("K3").Value = 71.60027500000001

# Row 4 (distance 1400)
$ws.Range("B4").Value = 130236.0
$ws.Range("C4").Value = 0.398
$ws.Range("E4").Value = 82.93398855001689
$ws.Range("H4").Value = 1.437473558644611
$ws.Range("I4").Value = 82.087
$ws.Range("J4").Value = 83.04
$ws.Range("K4").Value = 83.9138

# Row 5 (distance 1600)
$ws.Range("B5").Value = 145926.0
$ws.Range("C5").Value = 0.392
$ws.Range("E5").Value = 94.99810640461604
$ws.Range("H5").Value = 1.6673813243137179
$ws.Range("I5").Value = 94.0064
$ws.Range("J5").Value = 95.13
$ws.Range("K5").Value = 96.1203

# Row 6 (distance 1800)
$ws.Range("B6").Value = 106725.0
$ws.Range("C6").Value = 0.393
$ws.Range("E6").Value = 107.0268087917545
$ws.Range("H6").Value = 1.8443838565561554
$ws.Range("I6").Value = 105.9317
$ws.Range("J6").Value = 107.18
$ws.Range("K6").Value = 108.2829

# Row 7 (distance 2000)
$ws.Range("B7").Value = 121377.0
$ws.Range("C7").Value = 0.399
$ws.Range("E7").Value = 118.91322819562191
$ws.Range("H7").Value = 2.0537577044301254
$ws.Range("I7").Value = 117.6816
$ws.Range("J7").Value = 119.09
$ws.Range("K7").Value = 120.314

# Row 8 (distance 2200)
$ws.Range("B8").Value = 124758.0
$ws.Range("C8").Value = 0.423
$ws.Range("E8").Value = 130.81677295804678
$ws.Range("H8").Value = 2.264490349651805
$ws.Range("I8").Value = 129.4574
$ws.Range("J8").Value = 131.0
$ws.Range("K8").Value = 132.3659

# Row 9 (distance 2400)
$ws.Range("B9").Value = 92790.0
$ws.Range("C9").Value = 0.46
$ws.Range("E9").Value = 142.58305632719043
$ws.Range("H9").Value = 2.5254528284799864
$ws.Range("I9").Value = 141.0577
$ws.Range("J9").Value = 142.82
$ws.Range("K9").Value = 144.3197

# Row 10 (distance 2600)
$ws.Range("B10").Value = 14220.0
$ws.Range("C10").Value = 0.494
$ws.Range("E10").Value = 154.3668159774965
$ws.Range("H10").Value = 2.8989805076639117
$ws.Range("I10").Value = 152.6662
$ws.Range("J10").Value = 154.7
$ws.Range("K10").Value = 156.37052500000001
